# "green marked completed tests"
#
# The two "Wrong Answer in Penalty Box" bullet paragraphs get marked
# green -- the same accent6/BF-shaded green ("538135") already used for
# the other "completed" bullets throughout this document -- by adding:
#
#   <w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/>
#
# to the paragraph mark's run properties (w:pPr/w:rPr) and to the run
# properties (w:r/w:rPr) of each bullet's text run.

$d = $word.ActiveDocument

$colorRPr = '<w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr>'

$targets = @(
    "Verify that a wrong answer will send a player to the penalty box.",
    "Verify that getting an answer wrong will move on to the next player."
)

foreach ($targetText in $targets) {

    $found = $d.Content
    $found.Find.ClearFormatting()
    $ok = $found.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $targetText"
        continue
    }

    # Grab the enclosing paragraph's canonical OOXML so we can reuse its
    # existing attributes (paraId, rsids, pStyle/numPr, ...) verbatim.
    $para = $found.Paragraphs(1)
    $prng = $para.Range
    $pkgXml = $prng.WordOpenXML

    $pMatch = [regex]::Match($pkgXml, '<w:p\b.*?</w:p>')
    $pXml = $pMatch.Value

    # Give the paragraph mark (w:pPr/w:rPr) the green color.
    if ($pXml -match '</w:pPr>') {
        $newPXml = $pXml -replace '</w:pPr>', ($colorRPr + '</w:pPr>')
    } else {
        $newPXml = $pXml -replace '(<w:pPr>)', ('${1}' + $colorRPr)
    }

    # Give every run in the paragraph (w:r/w:rPr) the green color too.
    $newPXml = $newPXml -replace '(<w:r(?:\s[^>]*)?>)', ('${1}' + $colorRPr)

    $newPkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
        '<w:body>' + $newPXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # Replace just the matched run's text (NOT Paragraph.Range, which
    # would also grab the trailing paragraph mark and merge this
    # paragraph with the next one on Delete()).
    $found.Delete()
    $found.InsertXML($newPkg) | Out-Null

    Write-Output "Updated: $targetText"
}
